# Performed calibration of the needle: re-sort the curvature data rows in
# ascending order by the "time (s)" column (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies A2:D18 (row 1 is the header row).
$dataRange = $ws.Range("A2:D18")
$keyRange  = $ws.Range("A2:A18")

$dataRange.Sort($keyRange, 1)
